# timeliste.xlsx - "La til gitignore fil" edit
#
# The sheet already contains the blank template rows/styles for the
# timesheet grid (rows 8-14 were pre-formatted placeholder rows, row 1
# was a blank spacer row). This script fills in the new log entries that
# were added to the sheet and clears the now-unused spacer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 was just a blank spacer (styled empty cells, no real content) -
# clear it out so the sheet's used range starts at row 2.
$ws.Range("C1:F1").ClearContents()

# --- New timesheet entries -------------------------------------------------
# Row 8: Christian
$ws.Range("A8").Value = "Christian"
$ws.Range("B8").Value = "Kode html og css"
$ws.Range("C8").Value = 43742
$ws.Range("D8").Value = 0.375
$ws.Range("E8").Value = 0.666666666666667

# Row 9: Iselin
$ws.Range("A9").Value = "Iselin"
$ws.Range("B9").Value = "Går nøye igennom prosjektoppgaven"
$ws.Range("C9").Value = 43559
$ws.Range("D9").Value = 0.854166666666667
$ws.Range("E9").Value = 0.895833333333333

# Row 10: Iselin
$ws.Range("A10").Value = "Iselin"
$ws.Range("B10").Value = "Kladd nettstedskart, forslag til navn Cruise"
$ws.Range("C10").Value = 43562
$ws.Range("D10").Value = 0.875
$ws.Range("E10").Value = 0.958333333333333

# Row 11: Iselin
$ws.Range("A11").Value = "Iselin"
$ws.Range("B11").Value = "Kode html"
$ws.Range("C11").Value = 43564
$ws.Range("D11").Value = 0.854166666666667
$ws.Range("E11").Value = 1

# Row 12: Iselin
$ws.Range("A12").Value = "Iselin"
$ws.Range("B12").Value = "Kode html"
$ws.Range("C12").Value = 43565
$ws.Range("D12").Value = 0.916666666666667
$ws.Range("E12").Value = 0.979166666666667

# Row 13: Alle
$ws.Range("A13").Value = "Alle"
$ws.Range("B13").Value = "Chat-møte planlegging"
$ws.Range("C13").Value = 43565
$ws.Range("D13").Value = 0.895833333333333
$ws.Range("E13").Value = 0.9375

# Row 14: Eva
$ws.Range("A14").Value = "Eva"
$ws.Range("B14").Value = "Arbeid med underveisrapport"
$ws.Range("C14").Value = 43566
$ws.Range("D14").Value = 0.791666666666667
$ws.Range("E14").Value = 0.833333333333333

# (Column F on every row already holds the "=E-D" duration formula from
# the original template, so it recalculates automatically.)

# Iterative-calculation max-change tightened from 0.001 to 0.0001.
$excel.Iteration = $False
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# The author's cursor ended up on F14 after entering the new rows.
$ws.Range("F14").Select() | Out-Null
